{"js": "// Replace the 25 division-problem strings in the body's single table.\n// Each \"old\" string is unique in the document, so a plain search+replace\n// per pair is unambiguous and order-independent.\nconst replacements = [\n  [\"620\u00f78=\", \"541\u00f77=\"],\n  [\"616\u00f72=\", \"694\u00f74=\"],\n  [\"894\u00f77=\", \"680\u00f73=\"],\n  [\"705\u00f73=\", \"340\u00f77=\"],\n  [\"603\u00f77=\", \"228\u00f73=\"],\n  [\"425\u00f78=\", \"741\u00f79=\"],\n  [\"678\u00f79=\", \"825\u00f74=\"],\n  [\"212\u00f79=\", \"738\u00f73=\"],\n  [\"640\u00f75=\", \"463\u00f74=\"],\n  [\"421\u00f78=\", \"877\u00f72=\"],\n  [\"157\u00f79=\", \"571\u00f73=\"],\n  [\"146\u00f79=\", \"560\u00f74=\"],\n  [\"451\u00f75=\", \"711\u00f72=\"],\n  [\"288\u00f73=\", \"193\u00f75=\"],\n  [\"815\u00f76=\", \"663\u00f78=\"],\n  [\"408\u00f79=\", \"751\u00f74=\"],\n  [\"241\u00f73=\", \"498\u00f77=\"],\n  [\"528\u00f77=\", \"894\u00f77=\"],\n  [\"591\u00f78=\", \"432\u00f72=\"],\n  [\"356\u00f76=\", \"310\u00f74=\"],\n  [\"442\u00f77=\", \"668\u00f75=\"],\n  [\"593\u00f79=\", \"428\u00f78=\"],\n  [\"131\u00f79=\", \"717\u00f76=\"],\n  [\"571\u00f79=\", \"653\u00f74=\"],\n  [\"278\u00f74=\", \"287\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 division-problem strings in the document's single table.\n# Each \"old\" string is unique in the document, so Find/Replace per pair is\n# unambiguous and order-independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"620\u00f78=\", \"541\u00f77=\"),\n    @(\"616\u00f72=\", \"694\u00f74=\"),\n    @(\"894\u00f77=\", \"680\u00f73=\"),\n    @(\"705\u00f73=\", \"340\u00f77=\"),\n    @(\"603\u00f77=\", \"228\u00f73=\"),\n    @(\"425\u00f78=\", \"741\u00f79=\"),\n    @(\"678\u00f79=\", \"825\u00f74=\"),\n    @(\"212\u00f79=\", \"738\u00f73=\"),\n    @(\"640\u00f75=\", \"463\u00f74=\"),\n    @(\"421\u00f78=\", \"877\u00f72=\"),\n    @(\"157\u00f79=\", \"571\u00f73=\"),\n    @(\"146\u00f79=\", \"560\u00f74=\"),\n    @(\"451\u00f75=\", \"711\u00f72=\"),\n    @(\"288\u00f73=\", \"193\u00f75=\"),\n    @(\"815\u00f76=\", \"663\u00f78=\"),\n    @(\"408\u00f79=\", \"751\u00f74=\"),\n    @(\"241\u00f73=\", \"498\u00f77=\"),\n    @(\"528\u00f77=\", \"894\u00f77=\"),\n    @(\"591\u00f78=\", \"432\u00f72=\"),\n    @(\"356\u00f76=\", \"310\u00f74=\"),\n    @(\"442\u00f77=\", \"668\u00f75=\"),\n    @(\"593\u00f79=\", \"428\u00f78=\"),\n    @(\"131\u00f79=\", \"717\u00f76=\"),\n    @(\"571\u00f79=\", \"653\u00f74=\"),\n    @(\"278\u00f74=\", \"287\u00f79=\")\n)\n\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, $wdReplaceAll)\n}\n"}
